# Updated symbol list with refreshed Price (D) and Volume(1h) (E) figures.
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the workbook's existing inline-string cells)
# instead of auto-converting numeric-looking / percentage-looking text
# into Number/Percentage typed cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'303.33"
$ws.Range("E2").Value  = "'1.34%"
$ws.Range("D3").Value  = "'33.06"
$ws.Range("E3").Value  = "'5.63%"
$ws.Range("D4").Value  = "'4.937"
$ws.Range("E4").Value  = "'-3.25%"
$ws.Range("D5").Value  = "'0.07832"
$ws.Range("E5").Value  = "'-1.60%"
$ws.Range("D6").Value  = "'2.010"
$ws.Range("E6").Value  = "'-12.33%"
$ws.Range("D7").Value  = "'7.839"
$ws.Range("E7").Value  = "'0.91%"
$ws.Range("D8").Value  = "'3.809"
$ws.Range("E8").Value  = "'-1.45%"
$ws.Range("D9").Value  = "'0.9240"
$ws.Range("E9").Value  = "'0.18%"
$ws.Range("D10").Value = "'0.1754"
$ws.Range("E10").Value = "'1.17%"
$ws.Range("D11").Value = "'0.07814"
$ws.Range("E11").Value = "'3.58%"
$ws.Range("D12").Value = "'0.08679"
$ws.Range("E12").Value = "'-6.84%"
$ws.Range("D13").Value = "'0.03145"
$ws.Range("E13").Value = "'3.29%"
$ws.Range("D14").Value = "'0.1004"
$ws.Range("E14").Value = "'0.20%"
$ws.Range("D15").Value = "'0.001519"
$ws.Range("E15").Value = "'0.57%"
$ws.Range("D16").Value = "'0.005926"
$ws.Range("E16").Value = "'2.93%"
$ws.Range("E17").Value = "'-0.54%"
$ws.Range("E18").Value = "'-5.02%"
$ws.Range("D19").Value = "'0.3309"
$ws.Range("E19").Value = "'1.14%"
$ws.Range("D20").Value = "'0.1317"
$ws.Range("E20").Value = "'-1.29%"
$ws.Range("E21").Value = "'10.66%"
$ws.Range("D22").Value = "'0.1992"
$ws.Range("E22").Value = "'17.18%"
$ws.Range("D23").Value = "'0.04564"
$ws.Range("E23").Value = "'-1.41%"
$ws.Range("D24").Value = "'0.001226"
$ws.Range("E24").Value = "'-1.97%"
$ws.Range("D25").Value = "'0.004443"
$ws.Range("E25").Value = "'-0.79%"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'4.35%"
$ws.Range("D39").Value = "'0.01738"
$ws.Range("E39").Value = "'-1.35%"
$ws.Range("D40").Value = "'0.04793"
$ws.Range("E40").Value = "'3.44%"
$ws.Range("D41").Value = "'0.007541"
$ws.Range("E41").Value = "'8.20%"
$ws.Range("D42").Value = "'0.1362"
$ws.Range("E42").Value = "'0.03%"
$ws.Range("E43").Value = "'6.94%"
$ws.Range("D44").Value = "'0.01173"
$ws.Range("E44").Value = "'13.95%"
$ws.Range("D45").Value = "'0.00006233"
$ws.Range("E45").Value = "'-0.50%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.16%"
$ws.Range("D48").Value = "'0.8204"
$ws.Range("E48").Value = "'-29.06%"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.16%"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.16%"
